$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, mirroring the styling of existing headers (style 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows for columns I and J
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 6
